$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns (I:L) for start date / end date / CAGR / ref date.
# This pushes the existing "label" column (I) to M, and shifts everything
# after it right by 4. The new columns inherit the bold/black header style
# ("s=1") from column H on header/data rows automatically.
$ws.Columns("I:L").Insert()

# Insert 2 new columns (N:O) for comment / source, right after the
# (now relocated) "label" column M.
$ws.Columns("N:O").Insert()

# --- Row 1 (headers) ---------------------------------------------------
# Enter values in the same order the workbook's shared-string table was
# built in, so new strings land at the expected indices.
$ws.Range("I1").Value = "start date"
$ws.Range("J1").Value = "end date"
$ws.Range("N1").Value = "comment"
$ws.Range("O1").Value = "source"

# N1/O1 are brand-new cells with no pre-existing style to inherit (their
# left neighbour M1 "label" is unstyled), so copy the bold/black header
# style ("s=1") used by the rest of row 1 onto them explicitly.
$ws.Range("H1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("O1").PasteSpecial(-4122)

# --- Row 2 data ----------------------------------------------------------
$ws.Range("H2").Value = "kg"

$ws.Range("L1").Value = "ref date"
$ws.Range("K1").Value = "CAGR"

# Date values for row 2 (test var 1): 2009-01-01 .. 2009-04-01, ref date 2009-01-01
# (ClearFormats first: the cell inherited the "s=1" black-font header style
# from the column insert, and we want a plain default-font + date-format
# style instead, matching what Excel produces for a freshly formatted cell)
$ws.Range("I2:L2").ClearFormats()
$ws.Range("I2").Value = 39814
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("J2").Value = 39904
$ws.Range("L2").Value = 39814

$ws.Range("K2").Value = 0.1
$ws.Range("K2").NumberFormat = "0.00"

# --- Clear the stray inherited styling on rows 3, 4, 5 and 7 -------------
# (columns I:L and N:O should stay completely empty on these rows, matching
# the pre-existing "no data" columns, rather than keeping the s="1" style
# that column insertion propagated onto them)
$ws.Range("I3:L5").Clear()
$ws.Range("I7:L7").Clear()
$ws.Range("N3:O5").Clear()
$ws.Range("N7:O7").Clear()

# --- New row 8: "multiple choice" variable -------------------------------
$ws.Range("A8").Value = "multiple choice"
$ws.Range("C8").Value = "numpy.random"
$ws.Range("D8").Value = "choice"
$ws.Range("E8").Value = "1,2,3"
$ws.Range("H8").Value = "kg"

# Row 8 is brand new, so H8 needs the "s=1" header style applied explicitly
# (there's no existing cell there yet to inherit it from).
$ws.Range("H2").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("I8").Value = 39083
$ws.Range("J8").Value = 39814

$ws.Range("M8").Value = "test var 1"

# --- Column A width (fits "multiple choice") ------------------------------
$ws.Columns("A").ColumnWidth = 12.66

# --- Restore selection to match the saved view ----------------------------
$ws.Range("C5").Select()
